# Update "Generate Report for Handback" timestamps.
#
# Mapping of cells (by sheet) to their new timestamp values:
#   Overview!G2 : 2016-09-06 09:29:23 -> 2016-09-06 09:30:32
#   zh-cn!H2    : 2016-09-06 09:29:17 -> 2016-09-06 09:30:10
#   zh-cn!K2    : 2016-09-06 09:29:35 -> 2016-09-06 09:30:57
#   de-de!K2    : 2016-09-06 09:29:43 -> 2016-09-06 09:31:17

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 09:30:32"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 09:30:10"
$wsZhCn.Range("K2").Value = "2016-09-06 09:30:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-06 09:31:17"
